# Daily attendance processing - 2026-01-13 09:12:39
# Reverse the order of the comma-separated "Recorded By" entries in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCell = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162)
$lastRow = $lastCell.Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ","
        if ($parts.Count -gt 1) {
            for ($i = 0; $i -lt $parts.Count; $i++) {
                $parts[$i] = $parts[$i].Trim()
            }

            $n = $parts.Count
            $reversed = @()
            for ($i = $n - 1; $i -ge 0; $i--) {
                $reversed += $parts[$i]
            }

            $cell.Value = [string]::Join(", ", $reversed)
        }
    }
}
